$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.686.37"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "2.368.74"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'299.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").Value = "'97.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("D7").Value = "'0.560"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("D10").Value = "'33.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.55%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "'7.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.98%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "2.732.10"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "2.371.99"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "'0.815"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "45.650.93"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "'12.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.01%  "
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "'6.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'66.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'242.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("E24").Value = "  -6.15%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'1.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.26%  "
$ws.Range("D27").Value = "'38.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -13.58%  "
$ws.Range("E28").Value = "  -4.09%  "
$ws.Range("D29").Value = "'9.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("E30").Value = "  +15.53%  "
$ws.Range("D31").Value = "'20.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").Value = "'5.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.89%  "
$ws.Range("D34").Value = "'146.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'0.0763"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.76%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("D37").Value = "'1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.57%  "
$ws.Range("D38").Value = "'0.115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("D39").Value = "'15.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.38%  "
$ws.Range("D40").Value = "'3.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.08%  "
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Value = "'3.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("D43").Value = "1.941.14"
$ws.Range("E43").Value = "  +4.47%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'93.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -11.13%  "
$ws.Range("E47").Value = "  +6.62%  "
$ws.Range("D48").Value = "'98.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").Value = "'0.182"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.46%  "
$ws.Range("D50").Value = "2.602.35"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").Value = "'68.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.47%  "
